# Correction de l'évaluation des compétences
# Fixes six typo/grammar mistakes in the "comments" column (E) of the
# evaluation sheet, and moves the selection cursor to E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E4: "La documentation et le découpage ... a été ... effectué" -> plural agreement
$ws.Range("E4").Value = "La documentation et le découpage fonctionnel de notre projet ont été, à notre humble avis, très bien effectués. Cependant, la gestion de temps fut somme toute inefficace. Les durées associées aux cas d'utilisations étaient inadéquates. Il y a eu un manque de compréhension de notre part. Durant notre parcours collégial, nous étions habitués à planifier des durées par cas d'utilisations. Cependant, dans un contexte comme celui-ci, il aurait été plus pertinent de planifier par semaine."

# E9: "aurais" -> "aurait"
$ws.Range("E9").Value = "Nous aurions pu trouver un outil qui génère la documentation pour notre API, cela nous aurait permis d'avoir de la documentation programmeur directement lié au code qui nécessiterait moins de temps à maintenir."

# E12: observé -> observer, ces -> ses, intéragir -> interagir, confiant -> sûrs
$ws.Range("E12").Value = "Nous avons créé un document détaillé pour expliquer clairement l'utilisation du logiciel par n'importe quel type d'utilisateur. Lors de la démonstration de l'application à la cliente, nous avons pu observer l'un de ses utilisateurs finaux interagir avec notre application et ainsi connaître la profondeur et les aspects à couvrir dans notre document de formation. Nous sommes donc sûrs qu'il convient à ceux-ci.  "

# E14: vu -> vue, cour -> cours
$ws.Range("E14").Value = "Nous avons placé nos fichiers dans des fichiers qui ont du sens (exemple la documentation dans le fichier docs) avec une structure vue dans le cours de génie 1. Il est donc facile de se retrouver. Donc nous pensons donc que c'est l'archivage des documents est bien fait."

# E17: excéder -> excédé, prit -> pris, iréaliste -> irréaliste, remove extra "à" before "chaque semaine" (x2)
$ws.Range("E17").Value = "Comme écrit plus haut l'estimation du temps requis pour les cas d'utilisation n'a pas été très réaliste, ainsi beaucoup de tâches ont excédé le temps requis chaque semaine et d'autres ont pris beaucoup moins de temps chaque semaine. De plus, puisque l'estimation n'était pas réaliste, alors l'échéancier de travail a été très difficile à suivre à cause du temps qui était irréaliste."

# E18: suivis -> suivi, "quelque fois" -> "quelquefois", développement -> développements, aurrions -> aurions, dévelopement -> développement
$ws.Range("E18").Value = "Nous aurions pu mieux utiliser la méthodologie de travail `"AGILE`" pour mieux gérer nos cycles de développement. Par exemple, on aurait pu faire plus de rencontres AGILE et plus réviser le code de nos pairs pour produire un code plus homogène et simple à comprendre. En général, nous avons bien suivi nos devis de développement, mais quelquefois le temps pris pour faire certains développements était plus long que prévu, ainsi nous aurions pu mieux prévoir le développement de l'application."

# Move the active selection cursor to E4 (matches the saved view state in the workbook)
$ws.Range("E4").Select()
